$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 2438.75
$ws.Range("E2").Value = 2388.84

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 83.66068886991174
$ws.Range("E3").Value = 55.00000000000001

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 11.25441636945128
$ws.Range("E4").Value = 10

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 4.529395020177157
$ws.Range("E5").Value = 6.666666666666667

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 35.64
$ws.Range("E6").Value = 31

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 238.83
$ws.Range("E7").Value = 300

$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 879.5599999999999
$ws.Range("E8").Value = 868

$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.189213174595201
$ws.Range("E9").Value = 0.4444444444444444

$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 2.468562147318364
$ws.Range("E10").Value = 4.444444444444445

$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.5270340416268984
$ws.Range("E11").Value = 2.666666666666667

$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 1534.53
$ws.Range("E12").Value = 2000

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 3569.46
$ws.Range("E13").Value = 3510

$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 14.42
$ws.Range("E14").Value = 6.8

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 312.5
$ws.Range("E15").Value = 303

$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1.56
$ws.Range("E16").Value = 0.9

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 2.45
$ws.Range("E17").Value = 1

$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1.11
$ws.Range("E18").Value = 1.1

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 11.8
$ws.Range("E19").Value = 11.5

$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 18.67
$ws.Range("E20").Value = 2

$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 68.26000000000001
$ws.Range("E21").Value = 66.09999999999999

$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 3519.75
$ws.Range("E22").Value = 560

$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 4.16
$ws.Range("E23").Value = 0.7

$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 449.36
$ws.Range("E24").Value = 322

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1058.93
$ws.Range("E25").Value = 649

$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 9.6
$ws.Range("E26").Value = 8
